$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Recommandations")
$ws2 = $wb.Worksheets.Item("Top_YTD")

# --- Sheet "Recommandations": rewrite rows 2-44 (full re-sort + new rows 41-44) ---
$ws1.Range("A2").Value = 'CFAO MOTORS CI'
$ws1.Range("B2").Value = 0
$ws1.Range("C2").Value = 3
$ws1.Range("D2").Value = 2750
$ws1.Range("E2").Value = 920
$ws1.Range("F2").Value = '🟡 Observer'
$ws1.Range("G2").Value = '➖ Neutre'

$ws1.Range("A3").Value = 'BRVM - SERVICES PUBLICS'
$ws1.Range("B3").Value = 0
$ws1.Range("C3").Value = 6
$ws1.Range("D3").Value = 2490.43
$ws1.Range("E3").Value = 109.35
$ws1.Range("F3").Value = '🟡 Observer'
$ws1.Range("G3").Value = '➖ Neutre'

$ws1.Range("A4").Value = 'SETAO CI'
$ws1.Range("B4").Value = 0
$ws1.Range("C4").Value = 3
$ws1.Range("D4").Value = 2400
$ws1.Range("E4").Value = 875
$ws1.Range("F4").Value = '🟡 Observer'
$ws1.Range("G4").Value = '➖ Neutre'

$ws1.Range("A5").Value = 'AIR LIQUIDE CI'
$ws1.Range("B5").Value = 0
$ws1.Range("C5").Value = 3
$ws1.Range("D5").Value = 2040
$ws1.Range("E5").Value = 700
$ws1.Range("F5").Value = '🟡 Observer'
$ws1.Range("G5").Value = '➖ Neutre'

$ws1.Range("A6").Value = 'NEI-CEDA CI'
$ws1.Range("B6").Value = 0
$ws1.Range("C6").Value = 3
$ws1.Range("D6").Value = 1980
$ws1.Range("E6").Value = 665
$ws1.Range("F6").Value = '🟡 Observer'
$ws1.Range("G6").Value = '➖ Neutre'

$ws1.Range("A7").Value = 'UNIWAX CI'
$ws1.Range("B7").Value = 0
$ws1.Range("C7").Value = 2
$ws1.Range("D7").Value = 1870
$ws1.Range("E7").Value = 945
$ws1.Range("F7").Value = '🟡 Observer'
$ws1.Range("G7").Value = '➖ Neutre'

$ws1.Range("A8").Value = 'BRVM - AUTRES SECTEURS'
$ws1.Range("B8").Value = 0
$ws1.Range("C8").Value = 3
$ws1.Range("D8").Value = 1848.23
$ws1.Range("E8").Value = 611.46
$ws1.Range("F8").Value = '🟡 Observer'
$ws1.Range("G8").Value = '➖ Neutre'

$ws1.Range("A9").Value = 'BRVM - DISTRIBUTION'
$ws1.Range("B9").Value = 0
$ws1.Range("C9").Value = 3
$ws1.Range("D9").Value = 1186.15
$ws1.Range("E9").Value = 400.51
$ws1.Range("F9").Value = '🟡 Observer'
$ws1.Range("G9").Value = '➖ Neutre'

$ws1.Range("A10").Value = 'BRVM - TRANSPORT'
$ws1.Range("B10").Value = 0
$ws1.Range("C10").Value = 3
$ws1.Range("D10").Value = 1053.76
$ws1.Range("E10").Value = 356.15
$ws1.Range("F10").Value = '🟡 Observer'
$ws1.Range("G10").Value = '➖ Neutre'

$ws1.Range("A11").Value = 'BRVM - AGRICULTURE'
$ws1.Range("B11").Value = 0
$ws1.Range("C11").Value = 3
$ws1.Range("D11").Value = 1001.48
$ws1.Range("E11").Value = 337.87
$ws1.Range("F11").Value = '🟡 Observer'
$ws1.Range("G11").Value = '➖ Neutre'

$ws1.Range("A12").Value = 'BRVM - INDUSTRIE              (**)'
$ws1.Range("B12").Value = 0
$ws1.Range("C12").Value = 2
$ws1.Range("D12").Value = 461.08
$ws1.Range("E12").Value = 227.94
$ws1.Range("F12").Value = '🟡 Observer'
$ws1.Range("G12").Value = '➖ Neutre'

$ws1.Range("A13").Value = 'BRVM - INDUSTRIELS'
$ws1.Range("B13").Value = 0
$ws1.Range("C13").Value = 3
$ws1.Range("D13").Value = 422.4
$ws1.Range("E13").Value = 142.19
$ws1.Range("F13").Value = '🟡 Observer'
$ws1.Range("G13").Value = '➖ Neutre'

$ws1.Range("A14").Value = 'BRVM-PRESTIGE'
$ws1.Range("B14").Value = 0
$ws1.Range("C14").Value = 3
$ws1.Range("D14").Value = 404.81
$ws1.Range("E14").Value = 135.81
$ws1.Range("F14").Value = '🟡 Observer'
$ws1.Range("G14").Value = '➖ Neutre'

$ws1.Range("A15").Value = 'BRVM - FINANCES'
$ws1.Range("B15").Value = 0
$ws1.Range("C15").Value = 3
$ws1.Range("D15").Value = 395.35
$ws1.Range("E15").Value = 132.33
$ws1.Range("F15").Value = '🟡 Observer'
$ws1.Range("G15").Value = '➖ Neutre'

$ws1.Range("A16").Value = 'BRVM - SERVICES FINANCIERS'
$ws1.Range("B16").Value = 0
$ws1.Range("C16").Value = 3
$ws1.Range("D16").Value = 388.55
$ws1.Range("E16").Value = 130.05
$ws1.Range("F16").Value = '🟡 Observer'
$ws1.Range("G16").Value = '➖ Neutre'

$ws1.Range("A17").Value = 'BRVM - CONSOMMATION DISCRETIONNAIRE'
$ws1.Range("B17").Value = 0
$ws1.Range("C17").Value = 3
$ws1.Range("D17").Value = 370.87
$ws1.Range("E17").Value = 125.14
$ws1.Range("F17").Value = '🟡 Observer'
$ws1.Range("G17").Value = '➖ Neutre'

$ws1.Range("A18").Value = 'BRVM - ENERGIE'
$ws1.Range("B18").Value = 0
$ws1.Range("C18").Value = 3
$ws1.Range("D18").Value = 320.56
$ws1.Range("E18").Value = 106.42
$ws1.Range("F18").Value = '🟡 Observer'
$ws1.Range("G18").Value = '➖ Neutre'

$ws1.Range("A19").Value = 'BRVM - TELECOMMUNICATIONS'
$ws1.Range("B19").Value = 0
$ws1.Range("C19").Value = 3
$ws1.Range("D19").Value = 282.26
$ws1.Range("E19").Value = 94.45
$ws1.Range("F19").Value = '🟡 Observer'
$ws1.Range("G19").Value = '➖ Neutre'

$ws1.Range("A20").Value = 'BRVM - CONSOMMATION DE BASE          (**)'
$ws1.Range("B20").Value = 0
$ws1.Range("C20").Value = 1
$ws1.Range("D20").Value = 196.68
$ws1.Range("E20").Value = 196.68
$ws1.Range("F20").Value = '🟡 Observer'
$ws1.Range("G20").Value = '➖ Neutre'

$ws1.Range("A21").Value = 'BRVM - CONSOMMATION DE BASE             (**)'
$ws1.Range("B21").Value = 0
$ws1.Range("C21").Value = 1
$ws1.Range("D21").Value = 194.46
$ws1.Range("E21").Value = 194.46
$ws1.Range("F21").Value = '🟡 Observer'
$ws1.Range("G21").Value = '➖ Neutre'

$ws1.Range("A22").Value = 'BRVM-PRINCIPAL                (**)'
$ws1.Range("B22").Value = 0
$ws1.Range("C22").Value = 1
$ws1.Range("D22").Value = 192.95
$ws1.Range("E22").Value = 192.95
$ws1.Range("F22").Value = '🟡 Observer'
$ws1.Range("G22").Value = '➖ Neutre'

$ws1.Range("A23").Value = 'BRVM-PRINCIPAL                 (**)'
$ws1.Range("B23").Value = 0
$ws1.Range("C23").Value = 1
$ws1.Range("D23").Value = 192.29
$ws1.Range("E23").Value = 192.29
$ws1.Range("F23").Value = '🟡 Observer'
$ws1.Range("G23").Value = '➖ Neutre'

$ws1.Range("A24").Value = 'SETAO CI (STAC)'
$ws1.Range("B24").Value = 3
$ws1.Range("C24").Value = 0
$ws1.Range("D24").Value = 19.66
$ws1.Range("E24").Value = 7.36
$ws1.Range("F24").Value = '🟢 Achat'
$ws1.Range("G24").Value = '✅ Renforcer'

$ws1.Range("A25").Value = 'AIR LIQUIDE CI (SIVC)'
$ws1.Range("B25").Value = 2
$ws1.Range("C25").Value = 0
$ws1.Range("D25").Value = 12.42
$ws1.Range("E25").Value = 5.22
$ws1.Range("F25").Value = '🟡 Observer'
$ws1.Range("G25").Value = '➖ Neutre'

$ws1.Range("A26").Value = 'ECOBANK TRANS. INCORP. TG (ETIT)'
$ws1.Range("B26").Value = 2
$ws1.Range("C26").Value = 0
$ws1.Range("D26").Value = 11.76
$ws1.Range("E26").Value = 5.88
$ws1.Range("F26").Value = '🟡 Observer'
$ws1.Range("G26").Value = '➖ Neutre'

$ws1.Range("A27").Value = 'SICOR CI (SICC)'
$ws1.Range("B27").Value = 1
$ws1.Range("C27").Value = 0
$ws1.Range("D27").Value = 7.2
$ws1.Range("E27").Value = 7.2
$ws1.Range("F27").Value = '🟡 Observer'
$ws1.Range("G27").Value = '➖ Neutre'

$ws1.Range("A28").Value = 'SICABLE CI (CABC)'
$ws1.Range("B28").Value = 1
$ws1.Range("C28").Value = 0
$ws1.Range("D28").Value = 6.77
$ws1.Range("E28").Value = 6.77
$ws1.Range("F28").Value = '🟡 Observer'
$ws1.Range("G28").Value = '➖ Neutre'

$ws1.Range("A29").Value = 'SUCRIVOIRE (SCRC)'
$ws1.Range("B29").Value = 1
$ws1.Range("C29").Value = 0
$ws1.Range("D29").Value = 6.43
$ws1.Range("E29").Value = 6.43
$ws1.Range("F29").Value = '🟡 Observer'
$ws1.Range("G29").Value = '➖ Neutre'

$ws1.Range("A30").Value = 'SOCIETE IVOIRIENNE DE BANQUE  (SIBC)'
$ws1.Range("B30").Value = 1
$ws1.Range("C30").Value = 0
$ws1.Range("D30").Value = 4.95
$ws1.Range("E30").Value = 4.95
$ws1.Range("F30").Value = '🟡 Observer'
$ws1.Range("G30").Value = '➖ Neutre'

$ws1.Range("A31").Value = 'NSIA BANQUE COTE D''IVOIRE (NSBC)'
$ws1.Range("B31").Value = 1
$ws1.Range("C31").Value = 0
$ws1.Range("D31").Value = 3.14
$ws1.Range("E31").Value = 3.14
$ws1.Range("F31").Value = '🟡 Observer'
$ws1.Range("G31").Value = '➖ Neutre'

$ws1.Range("A32").Value = 'TOTAL'
$ws1.Range("B32").Value = 0
$ws1.Range("C32").Value = 3
$ws1.Range("D32").Value = 0
$ws1.Range("E32").Value = 0
$ws1.Range("F32").Value = '🟡 Observer'
$ws1.Range("G32").Value = '➖ Neutre'

$ws1.Range("A33").Value = 'SERVAIR ABIDJAN CI (ABJC)'
$ws1.Range("B33").Value = 1
$ws1.Range("C33").Value = 1
$ws1.Range("D33").Value = -0.11
$ws1.Range("E33").Value = 4.06
$ws1.Range("F33").Value = '🟡 Observer'
$ws1.Range("G33").Value = '👀 À surveiller'

$ws1.Range("A34").Value = 'FILTISAC CI (FTSC)'
$ws1.Range("B34").Value = 0
$ws1.Range("C34").Value = 1
$ws1.Range("D34").Value = -0.63
$ws1.Range("E34").Value = -0.63
$ws1.Range("F34").Value = '🟡 Observer'
$ws1.Range("G34").Value = '➖ Neutre'

$ws1.Range("A35").Value = 'BICI CI (BICC)'
$ws1.Range("B35").Value = 0
$ws1.Range("C35").Value = 1
$ws1.Range("D35").Value = -0.65
$ws1.Range("E35").Value = -0.65
$ws1.Range("F35").Value = '🟡 Observer'
$ws1.Range("G35").Value = '➖ Neutre'

$ws1.Range("A36").Value = 'CORIS BANK INTERNATIONAL (CBIBF)'
$ws1.Range("B36").Value = 0
$ws1.Range("C36").Value = 1
$ws1.Range("D36").Value = -1
$ws1.Range("E36").Value = -1
$ws1.Range("F36").Value = '🟡 Observer'
$ws1.Range("G36").Value = '➖ Neutre'

$ws1.Range("A37").Value = 'SOGB CI (SOGC)'
$ws1.Range("B37").Value = 0
$ws1.Range("C37").Value = 1
$ws1.Range("D37").Value = -1.33
$ws1.Range("E37").Value = -1.33
$ws1.Range("F37").Value = '🟡 Observer'
$ws1.Range("G37").Value = '➖ Neutre'

$ws1.Range("A38").Value = 'ECOBANK COTE D''''IVOIRE (ECOC)'
$ws1.Range("B38").Value = 0
$ws1.Range("C38").Value = 1
$ws1.Range("D38").Value = -1.54
$ws1.Range("E38").Value = -1.54
$ws1.Range("F38").Value = '🟡 Observer'
$ws1.Range("G38").Value = '➖ Neutre'

$ws1.Range("A39").Value = 'BERNABE CI (BNBC)'
$ws1.Range("B39").Value = 1
$ws1.Range("C39").Value = 2
$ws1.Range("D39").Value = -1.56
$ws1.Range("E39").Value = 6.79
$ws1.Range("F39").Value = '🟡 Observer'
$ws1.Range("G39").Value = '👀 À surveiller'

$ws1.Range("A40").Value = 'SAFCA CI (SAFC)'
$ws1.Range("B40").Value = 0
$ws1.Range("C40").Value = 1
$ws1.Range("D40").Value = -2.12
$ws1.Range("E40").Value = -2.12
$ws1.Range("F40").Value = '🟡 Observer'
$ws1.Range("G40").Value = '➖ Neutre'

$ws1.Range("A41").Value = 'ORANGE COTE D''IVOIRE (ORAC)'
$ws1.Range("B41").Value = 0
$ws1.Range("C41").Value = 1
$ws1.Range("D41").Value = -2.74
$ws1.Range("E41").Value = -2.74
$ws1.Range("F41").Value = '🟡 Observer'
$ws1.Range("G41").Value = '➖ Neutre'

$ws1.Range("A42").Value = 'UNIWAX CI (UNXC)'
$ws1.Range("B42").Value = 0
$ws1.Range("C42").Value = 1
$ws1.Range("D42").Value = -4.23
$ws1.Range("E42").Value = -4.23
$ws1.Range("F42").Value = '🟡 Observer'
$ws1.Range("G42").Value = '➖ Neutre'

$ws1.Range("A43").Value = 'SMB CI (SMBC)'
$ws1.Range("B43").Value = 0
$ws1.Range("C43").Value = 1
$ws1.Range("D43").Value = -5.11
$ws1.Range("E43").Value = -5.11
$ws1.Range("F43").Value = '🟡 Observer'
$ws1.Range("G43").Value = '➖ Neutre'

$ws1.Range("A44").Value = 'UNILEVER CI (UNLC)'
$ws1.Range("B44").Value = 0
$ws1.Range("C44").Value = 3
$ws1.Range("D44").Value = -22.47
$ws1.Range("E44").Value = -7.48
$ws1.Range("F44").Value = '🔴 Vente'
$ws1.Range("G44").Value = '⚠️ Risque de décrochage'

# --- Sheet "Top_YTD": rewrite rows 2-11 (re-sort + value updates) ---
$ws2.Range("A2").Value = 'BRVM - SERVICES PUBLICS'
$ws2.Range("B2").Value = 499443.39

$ws2.Range("A3").Value = 'CFAO MOTORS CI'
$ws2.Range("B3").Value = 104972.75

$ws2.Range("A4").Value = 'SETAO CI'
$ws2.Range("B4").Value = 72162.12

$ws2.Range("A5").Value = 'AIR LIQUIDE CI'
$ws2.Range("B5").Value = 47170

$ws2.Range("A6").Value = 'NEI-CEDA CI'
$ws2.Range("B6").Value = 43795.7

$ws2.Range("A7").Value = 'BRVM - AUTRES SECTEURS'
$ws2.Range("B7").Value = 36616.8

$ws2.Range("A8").Value = 'BRVM - DISTRIBUTION'
$ws2.Range("B8").Value = 12055.94

$ws2.Range("A9").Value = 'UNIWAX CI'
$ws2.Range("B9").Value = 10611.25

$ws2.Range("A10").Value = 'BRVM - TRANSPORT'
$ws2.Range("B10").Value = 9087.98

$ws2.Range("A11").Value = 'BRVM - AGRICULTURE'
$ws2.Range("B11").Value = 8064.33
